$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10, pushing the existing rows 10-36 down to 11-37.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the "new" weekly record.
# (Same market/category/etc. as its neighbours; only the date, volume and
# price columns differ from a generic template row.)
$ws.Range("A10").Value = 1
$ws.Range("B10").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C10").Value = "Arica y Parinacota"
$ws.Range("D10").Value = 44980
$ws.Range("E10").Value = 15
$ws.Range("F10").Value = 100112003
$ws.Range("G10").Value = "Ajo"
$ws.Range("H10").Value = "Chino"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 400
$ws.Range("K10").Value = 15000
$ws.Range("L10").Value = 16000
$ws.Range("M10").Value = 15500
$ws.Range("N10").Value = "$/caja 10 kilos"
$ws.Range("O10").Value = "China"
$ws.Range("P10").Value = 1550
$ws.Range("Q10").Value = 10
$ws.Range("R10").Value = "Hortaliza"
